# Helper: write a value into a range while forcing it to be stored as TEXT
# (matches the workbook's existing convention of storing percentage/ratio-style
# numbers as inline strings, e.g. "0.67", "82.54") and without leaving behind
# a stray custom number-format style on the cell.
function Set-TextValue($range, $value) {
    $savedStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $savedStyle
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" (Total) sheet: insert a new "2022-Q4" row right after the header,
#    pushing the existing 2022-Q3 / 2021-Q4 / 2020-Q4 rows down by one.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Shift existing data rows down (bottom-up so we don't clobber source data).
$total.Range("A4:D4").Copy($total.Range("A5:D5"))
$total.Range("A3:D3").Copy($total.Range("A4:D4"))
$total.Range("A2:D2").Copy($total.Range("A3:D3"))

# New first data row: 2022-Q4, same counts as the prior "latest" quarter.
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 8
$total.Range("D2").Value = 0.09

# Re-sequence the index column.
$total.Range("A2").Value = 0
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

# ---------------------------------------------------------------------------
# 2. Duplicate the current "2022-Q3" sheet so the untouched original data is
#    preserved under the "2022-Q3" name, placed right after the quarter
#    being edited.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item(2)
$q3.Copy($null, $q3)
$q3Copy = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------------
# 3. Turn the original sheet (still holding the old "2022-Q3" data) into the
#    new "2022-Q4" sheet by updating its holdings figures. Rename the
#    original first so the freshly made copy can reclaim the "2022-Q3" name.
# ---------------------------------------------------------------------------
$q3.Name = "2022-Q4"
$q3Copy.Name = "2022-Q3"

# Row 2 - 011097 达诚宜创精选混合A
Set-TextValue $q3.Range("D2") "0.67"
Set-TextValue $q3.Range("E2") "82.54"
Set-TextValue $q3.Range("F2") "3.31"
Set-TextValue $q3.Range("G2") "0.0222"
$q3.Range("H2").Value = 8

# Row 3 - 010301 达诚成长先锋混合A
Set-TextValue $q3.Range("D3") "0.51"
Set-TextValue $q3.Range("E3") "81.22"
Set-TextValue $q3.Range("F3") "3.22"
$q3.Range("H3").Value = 7

# Row 4 - 010808 达诚策略先锋混合A
Set-TextValue $q3.Range("D4") "0.32"
Set-TextValue $q3.Range("E4") "81.99"
Set-TextValue $q3.Range("F4") "3.25"
Set-TextValue $q3.Range("G4") "0.0104"
$q3.Range("H4").Value = 8

# Row 5 - 010809 达诚策略先锋混合C
Set-TextValue $q3.Range("E5") "81.99"
Set-TextValue $q3.Range("F5") "3.25"
Set-TextValue $q3.Range("G5") "0.0104"
$q3.Range("H5").Value = 8

# Row 6 - 010302 达诚成长先锋混合C
Set-TextValue $q3.Range("D6") "0.30"
Set-TextValue $q3.Range("E6") "81.22"
Set-TextValue $q3.Range("F6") "3.22"
Set-TextValue $q3.Range("G6") "0.0097"
$q3.Range("H6").Value = 7

# Row 7 - 011031 达诚价值先锋灵活配置混合C
Set-TextValue $q3.Range("D7") "0.22"
Set-TextValue $q3.Range("E7") "71.99"
Set-TextValue $q3.Range("F7") "3.41"
$q3.Range("H7").Value = 8

# Row 8 - 011030 达诚价值先锋灵活配置混合A
Set-TextValue $q3.Range("D8") "0.20"
Set-TextValue $q3.Range("E8") "71.99"
Set-TextValue $q3.Range("F8") "3.41"
Set-TextValue $q3.Range("G8") "0.0068"
$q3.Range("H8").Value = 8

# Row 9 - 011098 达诚宜创精选混合C
Set-TextValue $q3.Range("D9") "0.17"
Set-TextValue $q3.Range("E9") "82.54"
Set-TextValue $q3.Range("F9") "3.31"
Set-TextValue $q3.Range("G9") "0.0056"
$q3.Range("H9").Value = 8

# ---------------------------------------------------------------------------
# 4. Restore the original active-sheet selection ("总计").
# ---------------------------------------------------------------------------
$total.Activate()
